# TRDDimensions.xlsx - "Fix 13 and 14 row data"
#
# Rows 24/29 ("left" rows) and 28/33 ("right" rows) of the two module
# layers hold a pad-row count in column E that was wrong (14 / 13).
# The corrected value is 16 pads for all four rows. Column E is a raw
# input on both sheets; on Sheet1 the downstream Q/R/S/T columns are
# formulas and recalculate automatically, while on Sheet2 the same
# columns were pasted as static values, so they must be corrected by
# hand to match what the (now fixed) formulas on Sheet1 compute.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet1 - formulas recompute Q/R/S/T automatically once E changes
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("E24").Value = 16
$ws1.Range("E28").Value = 16
$ws1.Range("E29").Value = 16
$ws1.Range("E33").Value = 16

# ---------------------------------------------------------------
# Sheet2 - same raw input fix, plus the static (non-formula) copies
# of the dependent cells have to be updated explicitly
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("E24").Value = 16
$ws2.Range("Q24").Value = -362.1
$ws2.Range("S24").Value = 158.72

$ws2.Range("E28").Value = 16
$ws2.Range("R28").Value = 362.1
$ws2.Range("S28").Value = 158.72

$ws2.Range("E29").Value = 16
$ws2.Range("Q29").Value = -372.65999999999997
$ws2.Range("S29").Value = 161.76

$ws2.Range("E33").Value = 16
$ws2.Range("R33").Value = 372.65999999999997
$ws2.Range("S33").Value = 161.76

# ---------------------------------------------------------------
# Window/selection state, best-effort match of the recorded views
# ---------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A4:S33").Select() | Out-Null

$ws2.Activate()
$ws2.Range("L27").Select() | Out-Null

Write-Host "Fixed rows 24/28/29/33 (pad rows 13->16 / 14->16) on Sheet1 and Sheet2"
